$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.831039309501648
$ws.Range("B1").Value = 1.762451767921448
$ws.Range("C1").Value = 1.659937024116516
$ws.Range("D1").Value = 0.9971694946289062
$ws.Range("E1").Value = 0.6661605834960938
